$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H2").Value = 3.3
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 2.2
